# Apply the "PO Forecast" sheet addition and header renames
$wb = $excel.ActiveWorkbook

# --- 1. Rename header on "Weekly Quantity" sheet (B1) ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Rename header on "Monthly Trend" sheet (B1) ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the bold/centered header style from the Weekly Quantity sheet's A1
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-format style from the Weekly Quantity sheet's A2 onto column A (rows 2-64)
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A64").PasteSpecial(-4122)

# Forecast data rows 2-64: ds, PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(44990.99999999999,49,11.44015447076131,84.96956472350946),
    @(45095.99999999999,42,5.644236380584369,79.69191131980477),
    @(45116.99999999999,41,5.672631257475939,79.52457742243301),
    @(45123.99999999999,40,0.8579236144346419,77.47444897245343),
    @(45130.99999999999,40,5.904188725992722,78.01013502933863),
    @(45137.99999999999,39,0.8289256094657498,76.11353583274681),
    @(45144.99999999999,39,1.021712395875023,77.2817493776581),
    @(45158.99999999999,38,0.1805920384426122,75.80630920326614),
    @(45165.99999999999,38,-0.0007175177293758609,74.41642367226071),
    @(45172.99999999999,37,1.293728937648777,75.52051195653415),
    @(45179.99999999999,37,-1.848657107648902,72.81621042301073),
    @(45186.99999999999,36,-1.302027627304075,73.66834151825758),
    @(45193.99999999999,36,-2.145434579500737,73.61098113078675),
    @(45207.99999999999,35,-2.532673619617068,72.10149708316717),
    @(45214.99999999999,34,-2.562409059243838,73.3789964766293),
    @(45221.99999999999,34,-6.35820795052154,71.03938916656826),
    @(45228.99999999999,34,-3.650471017492221,70.61132990278101),
    @(45235.99999999999,33,-5.804989026081962,70.32748268630229),
    @(45242.99999999999,33,-3.033277723457284,70.14209911621225),
    @(45249.99999999999,32,-6.625108586250268,69.44158140163782),
    @(45256.99999999999,32,-4.429658801662914,67.59621552188469),
    @(45263.99999999999,31,-6.887645816146271,69.82751874763632),
    @(45270.99999999999,31,-6.293366012920988,67.24591258260251),
    @(45277.99999999999,30,-5.396095687701102,68.5240348548369),
    @(45298.99999999999,29,-6.175967931238206,65.82913574638786),
    @(45305.99999999999,29,-7.789464901477613,64.81214743714642),
    @(45312.99999999999,28,-7.229839260278166,66.21251400799085),
    @(45319.99999999999,28,-9.256494687251395,65.86887193033228),
    @(45340.99999999999,27,-12.75985489529726,64.07372591618518),
    @(45347.99999999999,26,-10.95237981297617,63.20199478934044),
    @(45361.99999999999,25,-13.13614830851853,62.85664137530733),
    @(45389.99999999999,23,-13.81414559125868,59.28407208662271),
    @(45396.99999999999,23,-15.94602312639076,61.5218944363608),
    @(45403.99999999999,23,-14.50225753787963,59.13457888820479),
    @(45410.99999999999,22,-18.35273342054244,61.52910370475438),
    @(45417.99999999999,22,-15.97863117989255,58.20677428112253),
    @(45424.99999999999,21,-16.84284402631207,55.77183885439798),
    @(45431.99999999999,21,-16.7136273589853,56.48932248597151),
    @(45438.99999999999,20,-16.28998789423687,58.62134258321237),
    @(45452.99999999999,19,-19.69541271987141,55.9418002783251),
    @(45473.99999999999,18,-20.36005265733079,56.00564287730968),
    @(45480.99999999999,18,-20.82896387389325,54.93274443256497),
    @(45487.99999999999,17,-22.68869760213064,55.14574275386525),
    @(45501.99999999999,16,-21.59612611980069,56.66324568531562),
    @(45508.99999999999,16,-21.09119241210526,53.59686535321292),
    @(45515.99999999999,15,-19.18729714472384,52.79049387636852),
    @(45522.99999999999,15,-24.50742736159432,50.54504557893179),
    @(45529.99999999999,15,-24.86874148817169,51.48720963289259),
    @(45536.99999999999,14,-22.35915007404475,52.34211695987424),
    @(45543.99999999999,14,-24.49103145888434,52.25494542645764),
    @(45550.99999999999,13,-23.49907724544061,49.27712741011707),
    @(45557.99999999999,13,-22.91458905622266,51.59991216806661),
    @(45564.99999999999,12,-26.6369035448695,51.9619989135785),
    @(45571.99999999999,12,-26.45282947889901,51.22325593272369),
    @(45578.99999999999,11,-24.01670138319271,46.71896252606422),
    @(45585.99999999999,11,-23.23146704726793,52.02293897469358),
    @(45592.99999999999,11,-24.01478178420362,50.83912563026286),
    @(45599.99999999999,10,-28.77705666290768,46.98415223804604),
    @(45606.99999999999,10,-26.91286572404256,47.97062562024498),
    @(45613.99999999999,9,-28.35309903641337,47.32231860709455),
    @(45620.99999999999,9,-27.48818421441822,46.81655847436143),
    @(45627.99999999999,8,-26.23998763013905,45.79599863684756),
    @(45634.99999999999,8,-27.93332168467216,43.54791712510469)
)
$numRows = $data.Count
$numCols = 4
$arr = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}
$wsForecast.Range("A2:D64").Value = $arr
